$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F (想去人数 / interested count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2215
$ws1.Range("F5").Value = 13253
$ws1.Range("F7").Value = 118
$ws1.Range("F10").Value = 1188
$ws1.Range("F12").Value = 13798
$ws1.Range("F13").Value = 14431
$ws1.Range("F21").Value = 40
$ws1.Range("F22").Value = 1100
$ws1.Range("F25").Value = 5481
$ws1.Range("F27").Value = 499
$ws1.Range("F28").Value = 340
$ws1.Range("F30").Value = 87

# Sheet "全部类型" (All Types) - column F (想去人数 / interested count) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2215
$ws4.Range("F5").Value = 13253
$ws4.Range("F8").Value = 118
$ws4.Range("F11").Value = 1188
$ws4.Range("F13").Value = 13798
$ws4.Range("F14").Value = 14431
$ws4.Range("F22").Value = 40
$ws4.Range("F23").Value = 1100
$ws4.Range("F26").Value = 5481
$ws4.Range("F28").Value = 499
$ws4.Range("F29").Value = 340
$ws4.Range("F31").Value = 87
